$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header / count row) tweaks ---
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# --- Row 2 (CON) tweaks ---
$ws.Range("B2").Value = 5.0175368920160865
$ws.Range("C2").ClearContents()
$ws.Range("E2").Value = 5.9453661633681243
# D2 (7.1594031091868082) is left untouched

# --- Row 3 (STR) tweaks ---
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 6.5446773553235111
$ws.Range("D3").Value = 6.1832668890764779
$ws.Range("E3").Value = 5.709867510217463

# --- Selection now only spans the touched block ---
$ws.Range("B1:E3").Select()
